$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.280.71"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.088.88"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "342.81"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5221"
$ws.Range("E7").Value = "  +1.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4399"
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.44"
$ws.Range("E9").Value = "  +3.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09336"
$ws.Range("E10").Value = "  +1.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.168"
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.73"
$ws.Range("E12").Value = "  -1.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.632"
$ws.Range("E13").Value = "  +4.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.892"
$ws.Range("E14").Value = "  +1.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.047.36"
$ws.Range("E15").Value = "  -2.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "101.22"
$ws.Range("E16").Value = "  +1.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001155"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.005"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.11"
$ws.Range("E19").Value = "  +1.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06674"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.327"
$ws.Range("E21").Value = "  +2.16%  "
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.289.59"
$ws.Range("E23").Value = "  +1.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.48"
$ws.Range("E24").Value = "  -1.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.292"
$ws.Range("E25").Value = "  -1.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.78"
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.31"
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.506"
$ws.Range("E28").Value = "  -1.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "132.89"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.657"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  +0.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.624"
$ws.Range("E34").Value = "  +10.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.869"
$ws.Range("E35").Value = "  -1.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.18"
$ws.Range("E36").Value = "  -2.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02625"
$ws.Range("E37").Value = "  +2.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06792"
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6971"
$ws.Range("E39").Value = "  +1.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.344"
$ws.Range("E40").Value = "  +3.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.49"
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2211"
$ws.Range("E42").Value = "  -1.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6792"
$ws.Range("E43").Value = "  +1.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.27"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.330"
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.002"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.372"
$ws.Range("E47").Value = "  +18.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.632"
$ws.Range("E48").Value = "  +0.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000345"
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("E50").Value = "  +8.63%  "
$ws.Range("E51").Value = "  -0.41%  "
